$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Act 1 description: replace the old placeholder story text with the
# real tutorial-stage blurb.
$ws.Range("B105").Value = "This is supposed to be an easy tutorial stage for the player to get used to the controls."

# Update the Act 2 description: replace the old placeholder text with the real
# blurb describing the "Arenas" combat encounters.
$ws.Range("B107").Value = 'Another easy tutorial stage. This one has some "Arenas" that you can''t pass until you beat all the enemies. More combat focused than the first level, but you can also skip a lot of encounters.'

# Insert two new rows for the new "Zero Limit Act 3" localization keys right
# after row 107 (pushing the old DEBUG_ENEMY rows down from 108-109 to 110-111).
$ws.Range("A108:A109").EntireRow.Insert()

$ws.Range("A108").Value = "ZERO_LIMIT_ACT_3"
$ws.Range("B108").Value = "Zero Limit Act 3"
$ws.Range("A109").Value = "ZERO_LIMIT_ACT_3_DESCRIPTION"
$ws.Range("B109").Value = "insert description here lol"

# Update the saved selection to match the author's final cursor position.
$ws.Range("B107").Select()
